$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the title text in B1: drop the period after "6.4.2.1"
$ws.Range("B1").Value = "6.4.2.1 Общий объем забора пресной воды "

# 2. Update data values for column L (year 2022)
$ws.Range("L5").Value = 8741.9

# L7 was a formula (=L5-L8); replace it with the plain cached result
$ws.Range("L7").Value = 8483.5

$ws.Range("L14").Value = 1327.6

$ws.Range("L18").Value = 54

# 3. Move the active selection from M4 to O2
$ws.Range("O2").Select()
